# Refresh the Fgf15-Fgfr3 LR-pair sheet with the new TPM-based NATMI output.
# The updated run adds "Inflammatory-Mac" as a sending cluster (alongside the
# existing "MuSCs") and now reports edges to all five target clusters
# (ECs, FAPs, MuSCs, Neutrophils, Resolving-Mac) for each sender, so the table
# grows from 2 senders x 3 targets (6 data rows) to 2 senders x 5 targets (10 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Inflammatory-Mac -> ECs
$ws.Cells.Item(2, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(2, 2).Value = "Fgf15"
$ws.Cells.Item(2, 3).Value = "Fgfr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.05008433333333333
$ws.Cells.Item(2, 8).Value = 0.150253
$ws.Cells.Item(2, 9).Value = 0.3054767171413236
$ws.Cells.Item(2, 10).Value = 0.3054767171413236
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.662797333333334
$ws.Cells.Item(2, 14).Value = 13.988392
$ws.Cells.Item(2, 15).Value = 0.7324994586787992
$ws.Cells.Item(2, 16).Value = 0.7324994586787993
$ws.Cells.Item(2, 17).Value = 0.2335330959084445
$ws.Cells.Item(2, 18).Value = 2.101797863176
$ws.Cells.Item(2, 19).Value = 0.2237615299449962
$ws.Cells.Item(2, 20).Value = 0.2237615299449962

# Row 3: Inflammatory-Mac -> FAPs
$ws.Cells.Item(3, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 2).Value = "Fgf15"
$ws.Cells.Item(3, 3).Value = "Fgfr3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.05008433333333333
$ws.Cells.Item(3, 8).Value = 0.150253
$ws.Cells.Item(3, 9).Value = 0.3054767171413236
$ws.Cells.Item(3, 10).Value = 0.3054767171413236
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.6655859999999999
$ws.Cells.Item(3, 14).Value = 1.996758
$ws.Cells.Item(3, 15).Value = 0.1045598489170565
$ws.Cells.Item(3, 16).Value = 0.1045598489170565
$ws.Cells.Item(3, 17).Value = 0.033335431086
$ws.Cells.Item(3, 18).Value = 0.300018879774
$ws.Cells.Item(3, 19).Value = 0.0319405993919752
$ws.Cells.Item(3, 20).Value = 0.0319405993919752

# Row 4: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(4, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 2).Value = "Fgf15"
$ws.Cells.Item(4, 3).Value = "Fgfr3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.05008433333333333
$ws.Cells.Item(4, 8).Value = 0.150253
$ws.Cells.Item(4, 9).Value = 0.3054767171413236
$ws.Cells.Item(4, 10).Value = 0.3054767171413236
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7894166666666665
$ws.Cells.Item(4, 14).Value = 2.36825
$ws.Cells.Item(4, 15).Value = 0.1240129561007488
$ws.Cells.Item(4, 16).Value = 0.1240129561007488
$ws.Cells.Item(4, 17).Value = 0.03953740747222222
$ws.Cells.Item(4, 18).Value = 0.35583666725
$ws.Cells.Item(4, 19).Value = 0.03788307071264783
$ws.Cells.Item(4, 20).Value = 0.03788307071264783

# Row 5: Inflammatory-Mac -> Neutrophils
$ws.Cells.Item(5, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(5, 2).Value = "Fgf15"
$ws.Cells.Item(5, 3).Value = "Fgfr3"
$ws.Cells.Item(5, 4).Value = "Neutrophils"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.05008433333333333
$ws.Cells.Item(5, 8).Value = 0.150253
$ws.Cells.Item(5, 9).Value = 0.3054767171413236
$ws.Cells.Item(5, 10).Value = 0.3054767171413236
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.08057833333333334
$ws.Cells.Item(5, 14).Value = 0.241735
$ws.Cells.Item(5, 15).Value = 0.01265840681643176
$ws.Cells.Item(5, 16).Value = 0.01265840681643176
$ws.Cells.Item(5, 17).Value = 0.004035712106111112
$ws.Cells.Item(5, 18).Value = 0.036321408955
$ws.Cells.Item(5, 19).Value = 0.003866848558522928
$ws.Cells.Item(5, 20).Value = 0.003866848558522928

# Row 6: Inflammatory-Mac -> Resolving-Mac
$ws.Cells.Item(6, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 2).Value = "Fgf15"
$ws.Cells.Item(6, 3).Value = "Fgfr3"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.05008433333333333
$ws.Cells.Item(6, 8).Value = 0.150253
$ws.Cells.Item(6, 9).Value = 0.3054767171413236
$ws.Cells.Item(6, 10).Value = 0.3054767171413236
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.16722
$ws.Cells.Item(6, 14).Value = 0.50166
$ws.Cells.Item(6, 15).Value = 0.02626932948696365
$ws.Cells.Item(6, 16).Value = 0.02626932948696365
$ws.Cells.Item(6, 17).Value = 0.008375102220000001
$ws.Cells.Item(6, 18).Value = 0.07537591998
$ws.Cells.Item(6, 19).Value = 0.008024668533181427
$ws.Cells.Item(6, 20).Value = 0.008024668533181427

# Row 7: MuSCs -> ECs
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Fgf15"
$ws.Cells.Item(7, 3).Value = "Fgfr3"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.1138703333333333
$ws.Cells.Item(7, 8).Value = 0.341611
$ws.Cells.Item(7, 9).Value = 0.6945232828586764
$ws.Cells.Item(7, 10).Value = 0.6945232828586764
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.662797333333334
$ws.Cells.Item(7, 14).Value = 13.988392
$ws.Cells.Item(7, 15).Value = 0.7324994586787992
$ws.Cells.Item(7, 16).Value = 0.7324994586787993
$ws.Cells.Item(7, 17).Value = 0.5309542866124445
$ws.Cells.Item(7, 18).Value = 4.778588579512
$ws.Cells.Item(7, 19).Value = 0.508737928733803
$ws.Cells.Item(7, 20).Value = 0.5087379287338031

# Row 8: MuSCs -> FAPs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Fgf15"
$ws.Cells.Item(8, 3).Value = "Fgfr3"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1138703333333333
$ws.Cells.Item(8, 8).Value = 0.341611
$ws.Cells.Item(8, 9).Value = 0.6945232828586764
$ws.Cells.Item(8, 10).Value = 0.6945232828586764
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.6655859999999999
$ws.Cells.Item(8, 14).Value = 1.996758
$ws.Cells.Item(8, 15).Value = 0.1045598489170565
$ws.Cells.Item(8, 16).Value = 0.1045598489170565
$ws.Cells.Item(8, 17).Value = 0.075790499682
$ws.Cells.Item(8, 18).Value = 0.682114497138
$ws.Cells.Item(8, 19).Value = 0.07261924952508128
$ws.Cells.Item(8, 20).Value = 0.07261924952508128

# Row 9: MuSCs -> MuSCs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Fgf15"
$ws.Cells.Item(9, 3).Value = "Fgfr3"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1138703333333333
$ws.Cells.Item(9, 8).Value = 0.341611
$ws.Cells.Item(9, 9).Value = 0.6945232828586764
$ws.Cells.Item(9, 10).Value = 0.6945232828586764
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.7894166666666665
$ws.Cells.Item(9, 14).Value = 2.36825
$ws.Cells.Item(9, 15).Value = 0.1240129561007488
$ws.Cells.Item(9, 16).Value = 0.1240129561007488
$ws.Cells.Item(9, 17).Value = 0.08989113897222221
$ws.Cells.Item(9, 18).Value = 0.80902025075
$ws.Cells.Item(9, 19).Value = 0.08612988538810099
$ws.Cells.Item(9, 20).Value = 0.086129885388101

# Row 10: MuSCs -> Neutrophils
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Fgf15"
$ws.Cells.Item(10, 3).Value = "Fgfr3"
$ws.Cells.Item(10, 4).Value = "Neutrophils"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.1138703333333333
$ws.Cells.Item(10, 8).Value = 0.341611
$ws.Cells.Item(10, 9).Value = 0.6945232828586764
$ws.Cells.Item(10, 10).Value = 0.6945232828586764
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.08057833333333334
$ws.Cells.Item(10, 14).Value = 0.241735
$ws.Cells.Item(10, 15).Value = 0.01265840681643176
$ws.Cells.Item(10, 16).Value = 0.01265840681643176
$ws.Cells.Item(10, 17).Value = 0.009175481676111112
$ws.Cells.Item(10, 18).Value = 0.082579335085
$ws.Cells.Item(10, 19).Value = 0.008791558257908833
$ws.Cells.Item(10, 20).Value = 0.008791558257908835

# Row 11: MuSCs -> Resolving-Mac
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Fgf15"
$ws.Cells.Item(11, 3).Value = "Fgfr3"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1138703333333333
$ws.Cells.Item(11, 8).Value = 0.341611
$ws.Cells.Item(11, 9).Value = 0.6945232828586764
$ws.Cells.Item(11, 10).Value = 0.6945232828586764
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.16722
$ws.Cells.Item(11, 14).Value = 0.50166
$ws.Cells.Item(11, 15).Value = 0.02626932948696365
$ws.Cells.Item(11, 16).Value = 0.02626932948696365
$ws.Cells.Item(11, 17).Value = 0.01904139714
$ws.Cells.Item(11, 18).Value = 0.17137257426
$ws.Cells.Item(11, 19).Value = 0.01824466095378222
$ws.Cells.Item(11, 20).Value = 0.01824466095378222

